# UrbanGulal_Daily_2026-01-13.xlsx update
# A new order ("Test3," / Square Heat Pad x1) came in at 2026-01-13 18:56.
# It is inserted as the new top data row (row 2) on "Daily Orders",
# pushing all the previous order rows down by one, and the "Summary"
# sheet's running totals are bumped accordingly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Daily Orders")

# Insert a new row right below the header, shifting existing orders down.
$ws.Rows.Item(2).Insert()

# Fill in the new order's details.
$ws.Range("A2").Value2 = 8
$ws.Range("B2").Value2 = "2026-01-13 18:56"
$ws.Range("C2").Value2 = "Sagar Borse"
# Phone number must stay text (leading-zero-safe), not be auto-numericised.
$ws.Range("D2").Value2 = "'7588930329"
$ws.Range("E2").Value2 = "Test3,"
$ws.Range("F2").Value2 = "Square Heat Pad x1"
$ws.Range("G2").Value2 = 50
$ws.Range("H2").Value2 = "NEW"
$ws.Range("I2").Value2 = "PENDING"

# Update the rollup counts/totals on the Summary sheet.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A2").Value2 = 8
$wsSummary.Range("B2").Value2 = 6
$wsSummary.Range("G2").Value2 = 375
